$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking text (European-style thousands
# separators, fixed decimal places like "138.00"). Assigning a plain string to
# .Value lets Excel auto-convert it to a Double (losing formatting / exactness),
# so force the cell to Text format first to preserve the literal string exactly,
# the same way pre-formatting a cell as Text before typing does in real Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D26", "D27", "D29", "D31", "D32", "D35", "D39", "D40", "D42", "D43", "D45", "D46", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Price (D column) updates
$ws.Range("D2").Value = "64.088.67"
$ws.Range("D3").Value = "2.760.33"
$ws.Range("D5").Value = "575.16"
$ws.Range("D6").Value = "158.89"
$ws.Range("D12").Value = "0.385"
$ws.Range("D13").Value = "3.249.81"
$ws.Range("D14").Value = "26.93"
$ws.Range("D15").Value = "63.734.89"
$ws.Range("D17").Value = "2.768.15"
$ws.Range("D18").Value = "12.17"
$ws.Range("D19").Value = "4.83"
$ws.Range("D20").Value = "360.04"
$ws.Range("D22").Value = "0.998"
$ws.Range("D26").Value = "8.52"
$ws.Range("D27").Value = "0.997"
$ws.Range("D29").Value = "7.36"
$ws.Range("D31").Value = "1.34"
$ws.Range("D32").Value = "168.89"
$ws.Range("D35").Value = "20.18"
$ws.Range("D39").Value = "349.64"
$ws.Range("D40").Value = "6.31"
$ws.Range("D42").Value = "39.10"
$ws.Range("D43").Value = "22.09"
$ws.Range("D45").Value = "0.0589"
$ws.Range("D46").Value = "138.00"
$ws.Range("D51").Value = "11.04"

# Volume(1h) (E column) updates
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("E10").Value = "  -13.58%  "
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("E21").Value = "  -6.10%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -8.65%  "
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("E25").Value = "  -4.04%  "
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -7.20%  "
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("E45").Value = "  -4.66%  "
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.24%  "
